$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "이웃집 백만장자(골드 리커버 에디션)"
$ws.Range("B3").Value = "토머스 J. 스탠리, 윌리엄 D. 댄코/홍정희 역"
$ws.Range("C3").Value = "리드리드출판"

$ws.Range("A4").Value = "부자 아빠 가난한 아빠 1(20주년 특별 기념판)"
$ws.Range("B4").Value = "로버트 기요사키/안진환 역"
$ws.Range("C4").Value = "민음인"

$ws.Range("A5").Value = "마음의 지혜"
$ws.Range("B5").Value = "김경일"
$ws.Range("C5").Value = "포레스트북스"

$ws.Range("A6").Value = "사장학개론"
$ws.Range("B6").Value = "김승호"
$ws.Range("C6").Value = "스노우폭스북스"

$ws.Range("A7").Value = "하늘과 바람과 별과 인간"
$ws.Range("B7").Value = "김상욱"
$ws.Range("C7").Value = "바다출판사"

$ws.Range("A8").Value = "알아차림에 대한 알아차림"
$ws.Range("B8").Value = "루퍼트 스파이라"
$ws.Range("C8").Value = "퍼블리온"

$ws.Range("A9").Value = "이미 늦었다고 생각하는 당신을 위한 김미경의 마흔 수업"
$ws.Range("B9").Value = "김미경"
$ws.Range("C9").Value = "어웨이크북스"

$ws.Range("A10").Value = "메리골드 마음 세탁소"
$ws.Range("B10").Value = "윤정은"
$ws.Range("C10").Value = "북로망스"

$ws.Range("A11").Value = "백만장자 메신저"
$ws.Range("B11").Value = "브렌든 버처드/위선주 역"
$ws.Range("C11").Value = "리더스북"

$ws.Range("A12").Value = "실전 매수매도 기법"
$ws.Range("B12").Value = "김영옥(데이짱)"
$ws.Range("C12").Value = "이레미디어"

$ws.Range("A13").Value = "부의 조건"
$ws.Range("B13").Value = "자유지성 아카데미 17인"
$ws.Range("C13").Value = "자유지성"

$ws.Range("A15").Value = "나의 돈 많은 고등학교 친구"
$ws.Range("B15").Value = "송희구"
$ws.Range("C15").Value = "서삼독"

$ws.Range("A16").Value = "2023 제14회 젊은작가상 수상작품집"
$ws.Range("B16").Value = "이미상, 김멜라, 성혜령, 이서수, 정선임"
$ws.Range("C16").Value = "문학동네"
